$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving things like leading/trailing
# zeros, multi-dot numbers ("61.467.70"), subscript digits, etc. Excel's COM
# .Value setter auto-coerces strings that look like plain numbers into real
# numeric cells, which would corrupt values such as "5.40" -> 5.4 or
# "0.997" -> 0.99699999999999999. Forcing the cell to Text format first
# prevents that coercion, and resetting the style back to Normal afterwards
# keeps the cell's formatting identical to the original (no explicit style).
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" '61.467.70'
Set-TextValue $ws "E2" '  -2.13%  '
Set-TextValue $ws "D3" '2.449.43'
Set-TextValue $ws "E3" '  -4.58%  '
Set-TextValue $ws "E4" '  -0.03%  '
Set-TextValue $ws "D5" '547.43'
Set-TextValue $ws "E5" '  -3.47%  '
Set-TextValue $ws "D6" '146.60'
Set-TextValue $ws "E6" '  -4.22%  '
Set-TextValue $ws "E7" '  +0.03%  '
Set-TextValue $ws "D8" '0.585'
Set-TextValue $ws "E8" '  -5.26%  '
Set-TextValue $ws "D9" '2.450.81'
Set-TextValue $ws "E9" '  -4.63%  '
Set-TextValue $ws "E10" '  -6.37%  '
Set-TextValue $ws "E11" '  -1.05%  '
Set-TextValue $ws "D12" '5.40'
Set-TextValue $ws "E12" '  -3.91%  '
Set-TextValue $ws "D13" '0.351'
Set-TextValue $ws "E13" '  -6.08%  '
Set-TextValue $ws "D14" '26.06'
Set-TextValue $ws "E14" '  -5.51%  '
Set-TextValue $ws "D15" '2.884.64'
Set-TextValue $ws "E15" '  -4.77%  '
Set-TextValue $ws "D16" '0.0000167'
Set-TextValue $ws "E16" '  -5.22%  '
Set-TextValue $ws "D17" '61.291.15'
Set-TextValue $ws "E17" '  -2.21%  '
Set-TextValue $ws "D18" '2.442.31'
Set-TextValue $ws "E18" '  -4.37%  '
Set-TextValue $ws "D19" '10.93'
Set-TextValue $ws "E19" '  -7.14%  '
Set-TextValue $ws "D20" '6.96'
Set-TextValue $ws "E20" '  -5.15%  '
Set-TextValue $ws "D21" '4.16'
Set-TextValue $ws "E21" '  -5.46%  '
Set-TextValue $ws "D22" '318.61'
Set-TextValue $ws "E23" '  +0.10%  '
Set-TextValue $ws "D24" '1.88'
Set-TextValue $ws "E24" '  +2.11%  '
Set-TextValue $ws "D25" '63.68'
Set-TextValue $ws "E25" '  -4.91%  '
Set-TextValue $ws "D26" '0.0₃0978'
Set-TextValue $ws "E26" '  -9.76%  '
Set-TextValue $ws "D27" '2.562.44'
Set-TextValue $ws "E27" '  -5.22%  '
Set-TextValue $ws "B28" 'Bittensor'
Set-TextValue $ws "C28" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws "D28" '538.54'
Set-TextValue $ws "E28" '  -2.85%  '
Set-TextValue $ws "B29" 'Binance-PegBSC-USD'
Set-TextValue $ws "C29" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws "D29" '0.997'
Set-TextValue $ws "E29" '  -0.31%  '
Set-TextValue $ws "D30" '1.47'
Set-TextValue $ws "E30" '  -6.31%  '
Set-TextValue $ws "B31" 'Aptos'
Set-TextValue $ws "C31" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws "D31" '7.76'
Set-TextValue $ws "E31" '  -2.44%  '
Set-TextValue $ws "B32" 'InternetComputer(DFINITY)'
Set-TextValue $ws "C32" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws "D32" '8.24'
Set-TextValue $ws "E32" '  -7.94%  '
Set-TextValue $ws "E33" '  -6.16%  '
Set-TextValue $ws "E34" '  -5.24%  '
Set-TextValue $ws "D35" '1.58'
Set-TextValue $ws "E35" '  -5.17%  '
Set-TextValue $ws "D36" '5.75'
Set-TextValue $ws "E36" '  -9.60%  '
Set-TextValue $ws "D37" '0.999'
Set-TextValue $ws "E37" '  +0.03%  '
Set-TextValue $ws "D38" '4.82'
Set-TextValue $ws "E38" '  -6.55%  '
Set-TextValue $ws "D39" '0.379'
Set-TextValue $ws "E39" '  -3.88%  '
Set-TextValue $ws "D40" '18.25'
Set-TextValue $ws "E40" '  -5.36%  '
Set-TextValue $ws "E41" '  -3.36%  '
Set-TextValue $ws "D42" '139.98'
Set-TextValue $ws "E42" '  -7.81%  '
Set-TextValue $ws "E43" '  +0.08%  '
Set-TextValue $ws "D44" '40.19'
Set-TextValue $ws "E44" '  -3.29%  '
Set-TextValue $ws "D45" '2.30'
Set-TextValue $ws "E45" '  -5.24%  '
Set-TextValue $ws "D46" '142.02'
Set-TextValue $ws "E46" '  -8.29%  '
Set-TextValue $ws "D47" '3.60'
Set-TextValue $ws "E47" '  -5.25%  '
Set-TextValue $ws "D48" '21.63'
Set-TextValue $ws "E48" '  -6.52%  '
Set-TextValue $ws "D49" '0.0534'
Set-TextValue $ws "E49" '  -6.25%  '
Set-TextValue $ws "D50" '0.589'
Set-TextValue $ws "E50" '  -4.81%  '
Set-TextValue $ws "D51" '0.0929'
Set-TextValue $ws "E51" '  -5.32%  '

Write-Host "Applied cryptos update"
